$wb = $excel.ActiveWorkbook

# Rows (by worksheet row number) and their new "want to go" counts (column F)
$updates = @{
    3  = 98
    5  = 11592
    6  = 803
    8  = 16
    12 = 22
    14 = 52
    17 = 330
    18 = 1340
    20 = 902
}

# Both "展览" and "全部类型" sheets contain the same event rows and need updating
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
